$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.7943844492440605
$ws.Range("F2").Value = 0.8479481641468682
$ws.Range("E3").Value = 0.7697624190064795
$ws.Range("E4").Value = 0.767170626349892
$ws.Range("F4").Value = 0.8241900647948164
$ws.Range("F5").Value = 0.8220302375809936
